# Add image scale and colormap parameters (I4/I2/I1/I3 low/high pixel
# limits, contrast, bias) as new columns J..Y on the "Sheet" worksheet.
# Headers go in row 1; the same constant value is written to every data
# row (2..34).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")

$headers = @(
    "I4 Low Pixel Limit",
    "I4 High Pixel Limit",
    "I2 Low Pixel Limit",
    "I2 High Pixel Limit",
    "I1 Low Pixel Limit",
    "I1 High Pixel Limit",
    "I3 Low Pixel Limit",
    "I3 High Pixel Limit",
    "I4 Contrast",
    "I4 Bias",
    "I2 Contrast",
    "I2 Bias",
    "I1 Contrast",
    "I1 Bias",
    "I3 Contrast",
    "I3 Bias"
)

$values = @(
    3.0954,
    11.9371,
    0.564228,
    5.61279,
    0.564228,
    5.44242,
    -1.58671,
    22.5305,
    3.47048,
    0.516276,
    1,
    0.5,
    1,
    0.5,
    2.37925,
    0.567057
)

# Columns J (10) through Y (25)
$startCol = 10
$lastRow = 34

for ($i = 0; $i -lt $headers.Count; $i++) {
    $col = $startCol + $i

    # Row 1 header
    $ws.Cells.Item(1, $col).Value = $headers[$i]

    # Data rows 2..34 all share the same value
    for ($r = 2; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, $col).Value = $values[$i]
    }
}
